$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" "45.450.62"
Set-TextValue "E2" "  +2.99%  "
Set-TextValue "D3" "2.428.00"
Set-TextValue "E3" "  -0.41%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  -0.07%  "
Set-TextValue "D5" "318.59"
Set-TextValue "E5" "  +3.45%  "
Set-TextValue "D6" "102.84"
Set-TextValue "E6" "  +4.58%  "
Set-TextValue "D7" "0.516"
Set-TextValue "E7" "  +0.74%  "
Set-TextValue "E8" "  -0.12%  "
Set-TextValue "D9" "0.530"
Set-TextValue "E9" "  +6.31%  "
Set-TextValue "D10" "35.62"
Set-TextValue "E10" "  +0.80%  "
Set-TextValue "D11" "0.0804"
Set-TextValue "E12" "  -2.18%  "
Set-TextValue "D13" "18.16"
Set-TextValue "E13" "  -2.74%  "
Set-TextValue "D14" "7.06"
Set-TextValue "E14" "  +1.72%  "
Set-TextValue "D15" "2.807.48"
Set-TextValue "E15" "  -0.11%  "
Set-TextValue "D16" "2.427.90"
Set-TextValue "E16" "  -0.33%  "
Set-TextValue "D17" "0.847"
Set-TextValue "E17" "  +1.56%  "
Set-TextValue "D18" "45.362.53"
Set-TextValue "E18" "  +2.83%  "
Set-TextValue "D19" "12.24"
Set-TextValue "E19" "  -0.52%  "
Set-TextValue "E20" "  -1.67%  "
Set-TextValue "E21" "  +1.78%  "
Set-TextValue "D22" "68.88"
Set-TextValue "E22" "  +0.49%  "
Set-TextValue "D23" "244.66"
Set-TextValue "E23" "  +1.99%  "
Set-TextValue "E24" "  -0.40%  "
Set-TextValue "E25" "  +0.59%  "
Set-TextValue "E26" "  -0.02%  "
Set-TextValue "D27" "25.68"
Set-TextValue "E27" "  +1.86%  "
Set-TextValue "E28" "  -1.18%  "
Set-TextValue "D29" "9.59"
Set-TextValue "E29" "  +1.21%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D30" "33.06"
Set-TextValue "E30" "  +0.83%  "
$ws.Range("B31").Value = "OKB"
$ws.Range("C31").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D31" "49.21"
Set-TextValue "E31" "  +2.67%  "
Set-TextValue "D32" "20.35"
Set-TextValue "E32" "  +9.44%  "
Set-TextValue "E33" "  +5.60%  "
Set-TextValue "E34" "  +0.85%  "
Set-TextValue "E35" "  +0.20%  "
Set-TextValue "E36" "  +0.84%  "
Set-TextValue "E37" "  -2.84%  "
Set-TextValue "E38" "  -0.90%  "
Set-TextValue "D39" "2.87"
Set-TextValue "E39" "  -2.43%  "
Set-TextValue "D40" "126.18"
Set-TextValue "E40" "  -4.19%  "
Set-TextValue "E41" "  -3.11%  "
Set-TextValue "E42" "  +0.61%  "
Set-TextValue "D43" "20.60"
Set-TextValue "E43" "  -3.25%  "
Set-TextValue "E44" "  +1.47%  "
Set-TextValue "D45" "1.926.57"
Set-TextValue "E45" "  -1.43%  "
Set-TextValue "E46" "  -2.73%  "
Set-TextValue "E47" "  +1.88%  "
Set-TextValue "D48" "1.79"
Set-TextValue "E48" "  +10.31%  "
Set-TextValue "E49" "  -2.30%  "
Set-TextValue "D50" "76.74"
Set-TextValue "E50" "  +4.70%  "
Set-TextValue "D51" "53.91"
